# ---------------------------------------------------------------------------
# Applies the "Log File.docx" edit described by the commit:
#   "Moved code from header to implementation files and refactored it so
#    it is easily readable"
#
# Net effect (per the supplied OOXML diff):
#   1. A large, purely-cosmetic clean-up: every <w:proofErr/> spell/grammar
#      marker is gone and runs that existed only because a proofErr split
#      them are merged back into single runs. The visible text is
#      unchanged by this part.
#   2. "09/1/2016 (2hr)" -> "09/1/2016 (5hr)"
#   3. The next-to-last paragraph's truncated sentence "...computer's
#      though" + "t process to the user." is completed to "...computer's
#      thought process to the user." and three new bullet paragraphs are
#      appended after it; the hidden "_GoBack" bookmark ends up collapsed
#      right after the final sentence, in the very last paragraph.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: strip <w:proofErr/> markers and merge the runs that were only
# split because of them. Word's "Save as XML" (WordOpenXML) serialization
# already performs exactly this normalization, so round-tripping the whole
# document through it gives us the clean structure for free, without
# touching any real content.
# ---------------------------------------------------------------------------
$cleanXml = $d.WordOpenXML
$d.Content.InsertXML($cleanXml) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: "09/1/2016 (2hr)" -> "09/1/2016 (5hr)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("09/1/2016 (2hr)", $false, $false, $false, $false, `
                         $false, $true, 1, $false, "09/1/2016 (5hr)", 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: fix up the last two paragraphs.
#
# Before this step the document ends with a single paragraph whose logical
# text (spread across a couple of runs, split by the hidden "_GoBack"
# bookmark) reads:
#   "Added functions in notification class to display computer's though" +
#   "t process to the user."
# i.e. "...computer's thought process to the user." once the bookmark
# split is ignored.
#
# We replace that whole paragraph (dropping the old bookmark in the
# process) with the completed sentence followed by three new bullet
# paragraphs, matching the style of their ListBullet neighbours, then we
# recreate the "_GoBack" bookmark collapsed right after the final period
# of the very last paragraph.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs($lastParaIndex)
$targetRange = $targetPara.Range

# Exclude the trailing paragraph mark from the replacement range so we can
# control paragraph breaks explicitly via embedded carriage returns.
$replaceRange = $d.Range($targetRange.Start, $targetRange.End - 1)

$cr = [char]13
$rsquo = [char]8217

$newText = "Added functions in notification class to display computer" + $rsquo + "s thought process to the user." + $cr + `
           "Moved the code properly in between header and implementation files." + $cr + `
           "Besides documentation and inlining some functions, the code is presentable." + $cr + `
           "Besides the help mode, everything else is implemented in the game as per specifications."

$replaceRange.Text = $newText

# ---------------------------------------------------------------------------
# Step 4: recreate the hidden "_GoBack" bookmark collapsed immediately
# after the last sentence of the final paragraph (matching the target
# structure). Adding a collapsed bookmark exactly at the document's very
# last character position is unreliable, so we temporarily append a
# placeholder character, anchor the bookmark before it, then remove the
# placeholder again.
# ---------------------------------------------------------------------------
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$endPos = $finalRange.End - 1

$placeholderRange = $d.Range($endPos, $endPos)
$placeholderRange.InsertAfter("X") | Out-Null

$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

$finalPara2 = $d.Paragraphs($d.Paragraphs.Count)
$finalRange2 = $finalPara2.Range
$delStart = $finalRange2.End - 2
$delEnd = $finalRange2.End - 1
$d.Range($delStart, $delEnd).Text = ""
